$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates reflecting the latest cryptocurrency price/volume scrape.
# Values that parse as plain numbers are written with a Text number format
# first so Excel keeps the exact original string (no trailing-zero loss,
# no float drift) instead of silently converting them to a Number.

$ws.Range("D2").Value = '24.599.68'
$ws.Range("E2").Value = '  +3.25%  '
$ws.Range("D3").Value = '1.696.54'
$ws.Range("E3").Value = '  +2.02%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.02'
$ws.Range("E5").Value = '  +2.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3943'
$ws.Range("E7").Value = '  +1.66%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4016'
$ws.Range("E8").Value = '  +1.16%  '
$ws.Range("E9").Value = '  +4.88%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.001'
$ws.Range("E10").Value = '  -0.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.71'
$ws.Range("E11").Value = '  +0.61%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08764'
$ws.Range("E12").Value = '  +1.09%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.221'
$ws.Range("E13").Value = '  +6.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.28'
$ws.Range("E14").Value = '  +2.75%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.154'
$ws.Range("E15").Value = '  +11.84%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001315'
$ws.Range("E16").Value = '  +0.69%  '
$ws.Range("D17").Value = '1.694.58'
$ws.Range("E17").Value = '  +1.81%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '99.71'
$ws.Range("E18").Value = '  +0.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07069'
$ws.Range("E19").Value = '  +2.75%  '
$ws.Range("E20").Value = '  +3.34%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.944'
$ws.Range("E21").Value = '  +4.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.40%  '
$ws.Range("E23").Value = '  +3.04%  '
$ws.Range("D24").Value = '24.604.39'
$ws.Range("E24").Value = '  +3.32%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.130'
$ws.Range("E25").Value = '  +10.44%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.336'
$ws.Range("E26").Value = '  +1.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.83'
$ws.Range("E27").Value = '  +5.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.69'
$ws.Range("E28").Value = '  +1.96%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '136.41'
$ws.Range("E29").Value = '  +4.86%  '
$ws.Range("E30").Value = '  +1.36%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.464'
$ws.Range("E31").Value = '  +9.71%  '
$ws.Range("D32").Value = '1.879.35'
$ws.Range("E32").Value = '  +1.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.088'
$ws.Range("E33").Value = '  -2.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08594'
$ws.Range("E34").Value = '  +0.94%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.147'
$ws.Range("E35").Value = '  +7.55%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '11.56'
$ws.Range("E36").Value = '  +10.69%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2739'
$ws.Range("E37").Value = '  +3.55%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.928'
$ws.Range("E38").Value = '  +0.39%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.49'
$ws.Range("E39").Value = '  +0.24%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09132'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02739'
$ws.Range("E41").Value = '  +8.83%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.481'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7668'
$ws.Range("E43").Value = '  +1.49%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7177'
$ws.Range("E44").Value = '  +1.92%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '15.63'
$ws.Range("E45").Value = '  +4.53%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.550'
$ws.Range("E46").Value = '  +5.70%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.220'
$ws.Range("E47").Value = '  +2.65%  '
$ws.Range("E48").Value = '  -0.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '140.91'
$ws.Range("E49").Value = '  +1.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.323'
$ws.Range("E50").Value = '  +8.63%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07988'
$ws.Range("E51").Value = '  +2.55%  '
